$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number format first for price cells that would otherwise be
# auto-converted to numeric values by Excel (single-decimal-point numeric strings).
$textCells = @("D5","D6","D16","D19","D20","D21","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the latest crypto data refresh.
$ws.Range('D2').Value = '60.113.00'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '2.570.48'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '505.36'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').Value = '151.93'
$ws.Range('E6').Value = '  -4.70%  '
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  -5.65%  '
$ws.Range('D9').Value = '2.574.31'
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('E10').Value = '  +7.34%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = '3.022.29'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '60.193.18'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '21.48'
$ws.Range('E16').Value = '  -2.68%  '
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '2.574.90'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -0.74%  '
$ws.Range('D20').Value = '343.57'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '10.37'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  -0.69%  '
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.684.23'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0837'
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '7.37'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '19.27'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '153.15'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '1.55'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '5.69'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '3.98'
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.18'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').Value = '0.845'
$ws.Range('E38').Value = '  +7.48%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '0.845'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.47'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '36.06'
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.73'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '296.97'
$ws.Range('E43').Value = '  -5.19%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.616'
$ws.Range('E44').Value = '  -2.91%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = '0.0991'
$ws.Range('E45').Value = '  -2.85%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0555'
$ws.Range('E47').Value = '  -3.40%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '19.65'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '4.83'
$ws.Range('E49').Value = '  -3.72%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0232'
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '10.31'
$ws.Range('E51').Value = '  +0.41%  '
